{"js": "// Update the division problems in the worksheet table.\n// The table's data rows (0, 4, 8, 12, 16) each hold 5 \"NN\u00f7N=\" problems\n// across columns 0-4; this lists the old -> new problem text for every\n// one of those 25 cells, in table (row, col) order.\nconst replacements = [\n  { row: 0, col: 0, from: \"36\u00f78=\", to: \"36\u00f74=\" },\n  { row: 0, col: 1, from: \"15\u00f72=\", to: \"21\u00f77=\" },\n  { row: 0, col: 2, from: \"51\u00f72=\", to: \"16\u00f75=\" },\n  { row: 0, col: 3, from: \"39\u00f75=\", to: \"13\u00f76=\" },\n  { row: 0, col: 4, from: \"76\u00f78=\", to: \"23\u00f72=\" },\n  { row: 4, col: 0, from: \"55\u00f77=\", to: \"61\u00f79=\" },\n  { row: 4, col: 1, from: \"26\u00f79=\", to: \"30\u00f74=\" },\n  { row: 4, col: 2, from: \"98\u00f72=\", to: \"62\u00f74=\" },\n  { row: 4, col: 3, from: \"72\u00f75=\", to: \"12\u00f73=\" },\n  { row: 4, col: 4, from: \"64\u00f74=\", to: \"14\u00f73=\" },\n  { row: 8, col: 0, from: \"98\u00f77=\", to: \"47\u00f73=\" },\n  { row: 8, col: 1, from: \"63\u00f75=\", to: \"58\u00f72=\" },\n  { row: 8, col: 2, from: \"60\u00f72=\", to: \"14\u00f76=\" },\n  { row: 8, col: 3, from: \"27\u00f75=\", to: \"27\u00f72=\" },\n  { row: 8, col: 4, from: \"22\u00f79=\", to: \"86\u00f76=\" },\n  { row: 12, col: 0, from: \"45\u00f79=\", to: \"91\u00f78=\" },\n  { row: 12, col: 1, from: \"10\u00f77=\", to: \"14\u00f79=\" },\n  { row: 12, col: 2, from: \"98\u00f73=\", to: \"73\u00f73=\" },\n  { row: 12, col: 3, from: \"23\u00f79=\", to: \"67\u00f77=\" },\n  { row: 12, col: 4, from: \"90\u00f73=\", to: \"94\u00f77=\" },\n  { row: 16, col: 0, from: \"81\u00f72=\", to: \"89\u00f74=\" },\n  { row: 16, col: 1, from: \"13\u00f78=\", to: \"61\u00f79=\" },\n  { row: 16, col: 2, from: \"90\u00f74=\", to: \"53\u00f72=\" },\n  { row: 16, col: 3, from: \"23\u00f79=\", to: \"51\u00f76=\" },\n  { row: 16, col: 4, from: \"93\u00f79=\", to: \"60\u00f74=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst currentValues = table.values;\n\nfor (const rep of replacements) {\n  const current = currentValues[rep.row][rep.col];\n  // Only touch cells that still hold the expected \"before\" text; if a cell\n  // already shows the target value there's nothing to do.\n  if (current === rep.from || current !== rep.to) {\n    const cell = table.getCell(rep.row, rep.col);\n    cell.value = rep.to;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division problems in the worksheet table.\n# Each data row of the table holds 5 \"NN\u00f7N=\" problems; this maps the\n# old problem text to the new problem text, by table (row, col) position\n# (1-based, matching the Word COM Table.Cell indexing).\n$replacements = @(\n  @{ Row = 1; Col = 1; From = \"36\u00f78=\"; To = \"36\u00f74=\" },\n  @{ Row = 1; Col = 2; From = \"15\u00f72=\"; To = \"21\u00f77=\" },\n  @{ Row = 1; Col = 3; From = \"51\u00f72=\"; To = \"16\u00f75=\" },\n  @{ Row = 1; Col = 4; From = \"39\u00f75=\"; To = \"13\u00f76=\" },\n  @{ Row = 1; Col = 5; From = \"76\u00f78=\"; To = \"23\u00f72=\" },\n  @{ Row = 5; Col = 1; From = \"55\u00f77=\"; To = \"61\u00f79=\" },\n  @{ Row = 5; Col = 2; From = \"26\u00f79=\"; To = \"30\u00f74=\" },\n  @{ Row = 5; Col = 3; From = \"98\u00f72=\"; To = \"62\u00f74=\" },\n  @{ Row = 5; Col = 4; From = \"72\u00f75=\"; To = \"12\u00f73=\" },\n  @{ Row = 5; Col = 5; From = \"64\u00f74=\"; To = \"14\u00f73=\" },\n  @{ Row = 9; Col = 1; From = \"98\u00f77=\"; To = \"47\u00f73=\" },\n  @{ Row = 9; Col = 2; From = \"63\u00f75=\"; To = \"58\u00f72=\" },\n  @{ Row = 9; Col = 3; From = \"60\u00f72=\"; To = \"14\u00f76=\" },\n  @{ Row = 9; Col = 4; From = \"27\u00f75=\"; To = \"27\u00f72=\" },\n  @{ Row = 9; Col = 5; From = \"22\u00f79=\"; To = \"86\u00f76=\" },\n  @{ Row = 13; Col = 1; From = \"45\u00f79=\"; To = \"91\u00f78=\" },\n  @{ Row = 13; Col = 2; From = \"10\u00f77=\"; To = \"14\u00f79=\" },\n  @{ Row = 13; Col = 3; From = \"98\u00f73=\"; To = \"73\u00f73=\" },\n  @{ Row = 13; Col = 4; From = \"23\u00f79=\"; To = \"67\u00f77=\" },\n  @{ Row = 13; Col = 5; From = \"90\u00f73=\"; To = \"94\u00f77=\" },\n  @{ Row = 17; Col = 1; From = \"81\u00f72=\"; To = \"89\u00f74=\" },\n  @{ Row = 17; Col = 2; From = \"13\u00f78=\"; To = \"61\u00f79=\" },\n  @{ Row = 17; Col = 3; From = \"90\u00f74=\"; To = \"53\u00f72=\" },\n  @{ Row = 17; Col = 4; From = \"23\u00f79=\"; To = \"51\u00f76=\" },\n  @{ Row = 17; Col = 5; From = \"93\u00f79=\"; To = \"60\u00f74=\" }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $r = $cell.Range\n    # Trim the trailing cell-mark/paragraph-mark characters before comparing.\n    $current = $r.Text.TrimEnd([char]13, [char]7)\n    if ($current -eq $rep.From -or $current -ne $rep.To) {\n        $r.Text = $rep.To\n    }\n}\n"}
